# Insert a new weekly price row for "Ciboulette" (Femacal de La Calera) at
# row 231, pushing the existing rows 231..301 down to 232..302.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(231).Insert()

$ws.Range("A231").Value = 3
$ws.Range("B231").Value = "Femacal de La Calera"
$ws.Range("C231").Value = "Coquimbo"
$ws.Range("D231").Value = 44663
$ws.Range("E231").Value = 5
$ws.Range("F231").Value = 100112039
$ws.Range("G231").Value = "Ciboulette"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 120
$ws.Range("K231").Value = 1500
$ws.Range("L231").Value = 1500
$ws.Range("M231").Value = 1500
$ws.Range("N231").Value = "$/docena de atados"
$ws.Range("O231").Value = "Provincia de Quillota"
$ws.Range("P231").Value = 500
$ws.Range("Q231").Value = 3
$ws.Range("R231").Value = "Hortaliza"
